# Update column G ("K" - strikeouts) values for rows 2-27 in kershaw_clayton.xlsx
# These values were regenerated from box-score "K" totals instead of the prior
# "Strike#" (pitch-level strike count) figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 6
    3  = 9
    4  = 4
    5  = 7
    6  = 10
    7  = 5
    8  = 8
    9  = 6
    10 = 4
    11 = 3
    12 = 6
    13 = 6
    14 = 10
    15 = 7
    16 = 4
    17 = 7
    18 = 4
    19 = 4
    20 = 2
    21 = 7
    22 = 3
    23 = 7
    24 = 13
    25 = 6
    26 = 3
    27 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
